$d = $word.ActiveDocument

# 1. Extend the closing sentence of the intro paragraph: the run that
#    holds only "." becomes " and reveals artworks upon reaching new
#    heights."
$d.Content.Find.Execute(
    "attaches itself to objects to pull itself to higher places.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "attaches itself to objects to pull itself to higher places and reveals artworks upon reaching new heights.",
    2) | Out-Null

# 2. Rewrite the Core Gameplay Loop bullet so it ends on "get higher and
#    reveal artworks." instead of "gather pellets and traverse
#    obstacles".
$d.Content.Find.Execute(
    "gather pellets and traverse obstacles",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "get higher and reveal artworks.",
    2) | Out-Null

# 3. Remove the "Walking, Jumping" bullet paragraph from the
#    Functionalities list entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Walking, Jumping`r") {
        $p.Range.Delete()
        break
    }
}
